# Update Backlog for Task
# Applies the changes made on the "Sprint 1" worksheet:
#  - Row 8 ("Design Admin dashboard page"): Status -> Finished, Day 2 -> 2, Day 3 -> 4
#  - Row 9 ("Design site buider page"): Status -> In Process, Day 3 -> 1
#  - Active selection moved to K12

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 1")

# Row 8: "Design Admin dashboard page"
$ws.Range("F8").Value = "Finished"
$ws.Range("I8").Value = 2
$ws.Range("J8").Value = 4

# Row 9: "Design site buider page"
$ws.Range("F9").Value = "In Process"
$ws.Range("J9").Value = 1

# Update the active cell / selection on the sheet
$ws.Activate()
$ws.Range("K12").Select()
